$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H18/H19/H20 values (J18/J19/J20 recalc automatically since they are formulas)
$ws.Range("H18").Value = 1100
$ws.Range("H19").Value = 1100
$ws.Range("H20").Value = 1100

# Update the frozen pane's top-left cell and the bottom-right pane selection
$ws.Range("A34").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B12").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A34").Select()
